# Update the Price (D) and Volume(1h) (E) columns for the refreshed crypto snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.648.07"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "'3.618.95"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'605.84"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'199.54"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.219"
$ws.Range("E9").Value = "  +9.68%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'53.62"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").Value = "'9.54"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'4.195.34"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "'677.06"
$ws.Range("E15").Value = "  +14.18%  "
$ws.Range("D16").Value = "'12.91"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "'70.744.30"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "'3.600.48"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "'18.98"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "'18.63"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").Value = "'5.36"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'104.86"
$ws.Range("D25").Value = "'4.62"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").Value = "'9.87"
$ws.Range("E28").Value = "  +4.40%  "
$ws.Range("D29").Value = "'34.16"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").Value = "'4.56"
$ws.Range("E30").Value = "  +6.72%  "
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").Value = "'12.16"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'63.26"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'3.949.74"
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("D36").Value = "'0.0₃0866"
$ws.Range("E36").Value = "  +7.40%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'3.03"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").Value = "'36.61"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").Value = "'495.33"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").Value = "'3.07"
$ws.Range("E44").Value = "  +9.89%  "
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").Value = "'3.44"
$ws.Range("E46").Value = "  +4.71%  "
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "'8.62"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'1.31"
$ws.Range("E51").Value = "  +2.17%  "
